$d = $word.ActiveDocument

# Prepend a literal "#" as its own run (matching the target formatting of
# the paragraph's existing first run) in front of the given paragraph's
# text, without merging it into the neighboring run. A plain
# Range.InsertBefore("#") would get folded into the adjacent run because
# its formatting is identical, so instead we insert a short-lived sibling
# paragraph holding "#" and then delete the paragraph mark that separates
# it from the original paragraph - this glues the two paragraphs back
# into one while keeping their runs distinct.
function Add-HashPrefix($idx) {
    $p = $d.Paragraphs($idx)
    $insertPoint = $p.Range.Duplicate
    $insertPoint.Collapse(1)
    $insertPoint.InsertParagraphBefore()

    $newPara = $d.Paragraphs($idx)
    $newPara.Range.Text = "#"

    $markRange = $newPara.Range
    $markRange.Collapse(0)
    $markRange.MoveEnd(1, 1)
    $markRange.Delete()
}

# Title, the "Pseudocode..." paragraph, and every Heading2 section title
# ("Module to Make Assignment Easier", "Part A".."Part E") each get a
# leading "#" run.
Add-HashPrefix 1
Add-HashPrefix 2
Add-HashPrefix 3
Add-HashPrefix 5
Add-HashPrefix 7
Add-HashPrefix 9
Add-HashPrefix 11
Add-HashPrefix 13
